$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing facility (id 7, "st. marks maternity hospital (smmh)")
# which is no longer present in the updated dataset.
$ws.Rows.Item(8).Delete()

# Update header text.
$ws.Range("B1").Value = "facility"

# Recode the facility names (also normalizes capitalization).
$ws.Range("B2").Value = "Central Hospital"
$ws.Range("B3").Value = "Military Hospital"
$ws.Range("B4").Value = "Missing"
$ws.Range("B5").Value = "Other"
$ws.Range("B6").Value = "Port Hospital"
$ws.Range("B7").Value = "St. Mark's Maternity Hospital (SMMH)"
